# Update market/profit figures (currentAveragePrice*, LevePrice*, LeveProfit*)
# across the Leve profit tracking sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# per the latest market data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 686.2857
$ws.Range("I18").Value = 686.2857
$ws.Range("K18").Value = 686.2857
$ws.Range("M18").Value = -402.2857

# Row 132
$ws.Range("H132").Value = 4267.4
$ws.Range("I132").Value = 2503.9778
$ws.Range("J132").Value = 20138.2
$ws.Range("K132").Value = 7511.9334
$ws.Range("L132").Value = 60414.60000000001
$ws.Range("M132").Value = -4981.9334
$ws.Range("N132").Value = -65474.60000000001

# Row 138
$ws.Range("H138").Value = 3155.3774
$ws.Range("J138").Value = 3675.838
$ws.Range("L138").Value = 11027.514
$ws.Range("N138").Value = -21307.514

# Row 141
$ws.Range("H141").Value = 11850.583
$ws.Range("I141").Value = 12645.637
$ws.Range("K141").Value = 37936.911
$ws.Range("M141").Value = -32756.911


$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 8160.0234
$ws.Range("I32").Value = 5282.3716
$ws.Range("K32").Value = 5282.3716
$ws.Range("M32").Value = -4995.3716

# Row 61
$ws.Range("H61").Value = 3587.0908
$ws.Range("I61").Value = 2582.25
$ws.Range("J61").Value = 6266.6665
$ws.Range("K61").Value = 2582.25
$ws.Range("L61").Value = 6266.6665
$ws.Range("M61").Value = -2370.25
$ws.Range("N61").Value = -6690.6665

# Row 74
$ws.Range("H74").Value = 129419.375
$ws.Range("I74").Value = 146765
$ws.Range("K74").Value = 146765
$ws.Range("M74").Value = -145891

# Row 77
$ws.Range("H77").Value = 129419.375
$ws.Range("I77").Value = 146765
$ws.Range("K77").Value = 733825
$ws.Range("M77").Value = -729457

# Row 102
$ws.Range("H102").Value = 5938.654
$ws.Range("I102").Value = 6677.85
$ws.Range("J102").Value = 3474.6667
$ws.Range("K102").Value = 6677.85
$ws.Range("L102").Value = 3474.6667
$ws.Range("M102").Value = -5055.85
$ws.Range("N102").Value = -6718.6667

# Row 136
$ws.Range("H136").Value = 3587.0908
$ws.Range("I136").Value = 2582.25
$ws.Range("J136").Value = 6266.6665
$ws.Range("K136").Value = 7746.75
$ws.Range("L136").Value = 18799.9995
$ws.Range("M136").Value = -5196.75
$ws.Range("N136").Value = -23899.9995


$ws = $wb.Worksheets.Item("BSM")
# Row 36
$ws.Range("H36").Value = 8984.5
$ws.Range("I36").Value = 8984.5
$ws.Range("K36").Value = 8984.5
$ws.Range("M36").Value = -8450.5


$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 128235.375
$ws.Range("I58").Value = 145911.86
$ws.Range("K58").Value = 145911.86
$ws.Range("M58").Value = -145708.86

# Row 59
$ws.Range("H59").Value = 105994.89
$ws.Range("I59").Value = 89500
$ws.Range("K59").Value = 89500
$ws.Range("M59").Value = -88355

# Row 62
$ws.Range("H62").Value = 3532.7778
$ws.Range("I62").Value = 3400
$ws.Range("K62").Value = 3400
$ws.Range("M62").Value = -2776

# Row 65
$ws.Range("H65").Value = 3532.7778
$ws.Range("I65").Value = 3400
$ws.Range("K65").Value = 17000
$ws.Range("M65").Value = -13880

# Row 68
$ws.Range("H68").Value = 27000
$ws.Range("I68").Value = 27000
$ws.Range("K68").Value = 27000
$ws.Range("M68").Value = -26251

# Row 71
$ws.Range("H71").Value = 27000
$ws.Range("I71").Value = 27000
$ws.Range("K71").Value = 81000
$ws.Range("M71").Value = -77256

# Row 134
$ws.Range("H134").Value = 26100.404
$ws.Range("I134").Value = 30804.828
$ws.Range("K134").Value = 92414.484
$ws.Range("M134").Value = -89879.484

# Row 136
$ws.Range("H136").Value = 128235.375
$ws.Range("I136").Value = 145911.86
$ws.Range("K136").Value = 437735.58
$ws.Range("M136").Value = -435185.58

# Row 141
$ws.Range("H141").Value = 460037.84
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 460037.84
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 460037.84
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -470397.84


$ws = $wb.Worksheets.Item("CUL")
# Row 16
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()

# Row 44
$ws.Range("H44").Value = 805.2
$ws.Range("I44").Value = 506.5
$ws.Range("J44").Value = 2000
$ws.Range("K44").Value = 1519.5
$ws.Range("L44").Value = 6000
$ws.Range("M44").Value = -1121.5
$ws.Range("N44").Value = -6796

# Row 59
$ws.Range("H59").Value = 3904.4
$ws.Range("J59").Value = 7999.5
$ws.Range("L59").Value = 23998.5
$ws.Range("N59").Value = -25078.5

# Row 60
$ws.Range("H60").Value = 277.16666
$ws.Range("I60").Value = 222.8
$ws.Range("J60").Value = 549
$ws.Range("K60").Value = 668.4000000000001
$ws.Range("L60").Value = 1647
$ws.Range("M60").Value = -417.4000000000001
$ws.Range("N60").Value = -2149

# Row 121
$ws.Range("H121").Value = 712.0769
$ws.Range("J121").Value = 917.44446
$ws.Range("L121").Value = 2752.33338
$ws.Range("N121").Value = -5372.33338


$ws = $wb.Worksheets.Item("GSM")
# Row 39
$ws.Range("H39").Value = 49989
$ws.Range("J39").Value = 49989
$ws.Range("L39").Value = 49989
$ws.Range("N39").Value = -51053

# Row 80
$ws.Range("H80").Value = 3550
$ws.Range("I80").Value = 3059.5833
$ws.Range("K80").Value = 3059.5833
$ws.Range("M80").Value = -2061.5833

# Row 83
$ws.Range("H83").Value = 3550
$ws.Range("I83").Value = 3059.5833
$ws.Range("K83").Value = 15297.9165
$ws.Range("M83").Value = -10305.9165

# Row 126
$ws.Range("H126").Value = 6007.815
$ws.Range("I126").Value = 5904.5
$ws.Range("J126").Value = 6158.091
$ws.Range("K126").Value = 17713.5
$ws.Range("L126").Value = 18474.273
$ws.Range("M126").Value = -15243.5
$ws.Range("N126").Value = -23414.273

# Row 135
$ws.Range("H135").Value = 40937.5
$ws.Range("I135").Value = 30000
$ws.Range("K135").Value = 30000
$ws.Range("M135").Value = -24930


$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 52089.363
$ws.Range("I22").Value = 86043.766
$ws.Range("K22").Value = 86043.766
$ws.Range("M22").Value = -85748.766

# Row 27
$ws.Range("H27").Value = 52089.363
$ws.Range("I27").Value = 86043.766
$ws.Range("K27").Value = 86043.766
$ws.Range("M27").Value = -85936.766

# Row 68
$ws.Range("H68").Value = 2961
$ws.Range("J68").Value = 4999
$ws.Range("L68").Value = 4999
$ws.Range("N68").Value = -6497

# Row 71
$ws.Range("H71").Value = 2961
$ws.Range("J71").Value = 4999
$ws.Range("L71").Value = 24995
$ws.Range("N71").Value = -32483

# Row 82
$ws.Range("H82").Value = 3018.6667
$ws.Range("I82").Value = 1861.8334
$ws.Range("J82").Value = 3404.2778
$ws.Range("K82").Value = 1861.8334
$ws.Range("L82").Value = 3404.2778
$ws.Range("M82").Value = -1500.8334
$ws.Range("N82").Value = -4126.2778

# Row 85
$ws.Range("H85").Value = 3018.6667
$ws.Range("I85").Value = 1861.8334
$ws.Range("J85").Value = 3404.2778
$ws.Range("K85").Value = 1861.8334
$ws.Range("L85").Value = 3404.2778
$ws.Range("M85").Value = -613.8334
$ws.Range("N85").Value = -5900.2778

# Row 132
$ws.Range("H132").Value = 33318.9
$ws.Range("I132").Value = 40904.418
$ws.Range("K132").Value = 122713.254
$ws.Range("M132").Value = -120183.254

# Row 135
$ws.Range("H135").Value = 79997.5
$ws.Range("J135").Value = 79997.5
$ws.Range("L135").Value = 79997.5
$ws.Range("N135").Value = -90137.5


$ws = $wb.Worksheets.Item("WVR")
# Row 41
$ws.Range("H41").Value = 27552.715
$ws.Range("J41").Value = 26907
$ws.Range("L41").Value = 26907
$ws.Range("N41").Value = -27687

# Row 113
$ws.Range("H113").Value = 935.24
$ws.Range("I113").Value = 453.93332
$ws.Range("J113").Value = 1657.2
$ws.Range("K113").Value = 1361.79996
$ws.Range("L113").Value = 4971.6
$ws.Range("M113").Value = 808.2000400000002
$ws.Range("N113").Value = -9311.6

# Row 126
$ws.Range("H126").Value = 85490.086
$ws.Range("I126").Value = 85490.086
$ws.Range("K126").Value = 256470.258
$ws.Range("M126").Value = -254000.258
